# "download articles with pandoc title blocks"
#
# The title paragraph ("Day After Day - July-August 1934", styled Heading1)
# and the byline paragraph ("By Dorothy Day", manually bolded) are rebuilt
# as a pandoc-style title block:
#   - paragraph 1 becomes style "Title", text split into one run per
#     word/punctuation token (pandoc's tokenisation of the title string)
#   - paragraph 2 becomes style "Authors", holding just "Dorothy Day"
#     (the "By " prefix and manual bold formatting are dropped), again
#     split into one run per token
#
# Both target paragraph styles ("Title" / "Authors") already exist in
# styles.xml, so we only need to point w:pStyle at them.

$d = $word.ActiveDocument

$wordNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function New-RunXml([string]$text) {
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    return "<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
}

# --- Paragraph 1: title -----------------------------------------------
$titleTokens = @("Day", " ", "After", " ", "Day", " ", "-", " ", "July", "-", "August", " ", "1934")
$titleRuns = ($titleTokens | ForEach-Object { New-RunXml $_ }) -join ""
$titleXml = "<w:p xmlns:w=`"$wordNs`"><w:pPr><w:pStyle w:val=`"Title`"/></w:pPr>$titleRuns</w:p>"

$p1 = $d.Paragraphs.Item(1)
[void]$p1.Range.InsertXML($titleXml)

# --- Paragraph 2: authors ----------------------------------------------
$authorTokens = @("Dorothy", " ", "Day")
$authorRuns = ($authorTokens | ForEach-Object { New-RunXml $_ }) -join ""
$authorsXml = "<w:p xmlns:w=`"$wordNs`"><w:pPr><w:pStyle w:val=`"Authors`"/></w:pPr>$authorRuns</w:p>"

$p2 = $d.Paragraphs.Item(2)
[void]$p2.Range.InsertXML($authorsXml)

Write-Output "Title paragraph: [$($d.Paragraphs.Item(1).Range.Text)]"
Write-Output "Authors paragraph: [$($d.Paragraphs.Item(2).Range.Text)]"
